# Applies the "Updated cryptos list" price/volume refresh to the crypto table.
# Each entry is (cell, new text value). Cells in column D whose new value looks
# like a plain number are pre-formatted as Text so Excel keeps them as strings
# (matching the inline-string storage used throughout the sheet) instead of
# silently converting them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '25.995.29'; ForceText = $False }
    @{ Cell = 'E2'; Value = '  +0.09%  '; ForceText = $False }
    @{ Cell = 'D3'; Value = '1.633.50'; ForceText = $False }
    @{ Cell = 'E3'; Value = '  -0.44%  '; ForceText = $False }
    @{ Cell = 'E4'; Value = '  +0.20%  '; ForceText = $False }
    @{ Cell = 'E5'; Value = '  -0.96%  '; ForceText = $False }
    @{ Cell = 'E6'; Value = '  -0.52%  '; ForceText = $False }
    @{ Cell = 'E7'; Value = '  +0.20%  '; ForceText = $False }
    @{ Cell = 'E8'; Value = '  -1.99%  '; ForceText = $False }
    @{ Cell = 'D9'; Value = '0.0625'; ForceText = $True }
    @{ Cell = 'E9'; Value = '  -2.12%  '; ForceText = $False }
    @{ Cell = 'D10'; Value = '18.53'; ForceText = $True }
    @{ Cell = 'E10'; Value = '  -5.67%  '; ForceText = $False }
    @{ Cell = 'D11'; Value = '0.0792'; ForceText = $True }
    @{ Cell = 'E11'; Value = '  -0.32%  '; ForceText = $False }
    @{ Cell = 'D12'; Value = '1.860.54'; ForceText = $False }
    @{ Cell = 'E12'; Value = '  -0.45%  '; ForceText = $False }
    @{ Cell = 'E13'; Value = '  -2.05%  '; ForceText = $False }
    @{ Cell = 'D14'; Value = '1.620.35'; ForceText = $False }
    @{ Cell = 'E14'; Value = '  -1.34%  '; ForceText = $False }
    @{ Cell = 'E15'; Value = '  -2.83%  '; ForceText = $False }
    @{ Cell = 'D16'; Value = '25.993.29'; ForceText = $False }
    @{ Cell = 'E17'; Value = '  -2.64%  '; ForceText = $False }
    @{ Cell = 'D18'; Value = '61.82'; ForceText = $True }
    @{ Cell = 'E18'; Value = '  -1.69%  '; ForceText = $False }
    @{ Cell = 'D20'; Value = '190.28'; ForceText = $True }
    @{ Cell = 'E20'; Value = '  -1.41%  '; ForceText = $False }
    @{ Cell = 'D21'; Value = '4.25'; ForceText = $True }
    @{ Cell = 'E21'; Value = '  -2.55%  '; ForceText = $False }
    @{ Cell = 'E22'; Value = '  -3.71%  '; ForceText = $False }
    @{ Cell = 'E23'; Value = '  -2.07%  '; ForceText = $False }
    @{ Cell = 'D24'; Value = '0.134'; ForceText = $True }
    @{ Cell = 'E24'; Value = '  +0.42%  '; ForceText = $False }
    @{ Cell = 'D25'; Value = '143.28'; ForceText = $True }
    @{ Cell = 'E25'; Value = '  -0.87%  '; ForceText = $False }
    @{ Cell = 'E26'; Value = '  +0.26%  '; ForceText = $False }
    @{ Cell = 'E27'; Value = '  -2.63%  '; ForceText = $False }
    @{ Cell = 'E28'; Value = '  -2.36%  '; ForceText = $False }
    @{ Cell = 'D29'; Value = '15.20'; ForceText = $True }
    @{ Cell = 'E29'; Value = '  -2.06%  '; ForceText = $False }
    @{ Cell = 'E30'; Value = '  -1.22%  '; ForceText = $False }
    @{ Cell = 'E31'; Value = '  -2.90%  '; ForceText = $False }
    @{ Cell = 'E32'; Value = '  -3.01%  '; ForceText = $False }
    @{ Cell = 'E33'; Value = '  -4.04%  '; ForceText = $False }
    @{ Cell = 'E34'; Value = '  -1.50%  '; ForceText = $False }
    @{ Cell = 'D36'; Value = '0.872'; ForceText = $True }
    @{ Cell = 'E36'; Value = '  -3.63%  '; ForceText = $False }
    @{ Cell = 'D37'; Value = '1.136.00'; ForceText = $False }
    @{ Cell = 'E37'; Value = '  +0.13%  '; ForceText = $False }
    @{ Cell = 'E38'; Value = '  -1.38%  '; ForceText = $False }
    @{ Cell = 'D39'; Value = '0.524'; ForceText = $True }
    @{ Cell = 'E39'; Value = '  -3.31%  '; ForceText = $False }
    @{ Cell = 'E40'; Value = '  -1.36%  '; ForceText = $False }
    @{ Cell = 'D41'; Value = '98.71'; ForceText = $True }
    @{ Cell = 'E41'; Value = '  -0.75%  '; ForceText = $False }
    @{ Cell = 'E42'; Value = '  -1.87%  '; ForceText = $False }
    @{ Cell = 'E43'; Value = '  -4.64%  '; ForceText = $False }
    @{ Cell = 'D44'; Value = '1.770.94'; ForceText = $False }
    @{ Cell = 'E44'; Value = '  -0.49%  '; ForceText = $False }
    @{ Cell = 'E45'; Value = '  -0.63%  '; ForceText = $False }
    @{ Cell = 'D46'; Value = '55.15'; ForceText = $True }
    @{ Cell = 'E46'; Value = '  -2.66%  '; ForceText = $False }
    @{ Cell = 'D47'; Value = '0.0528'; ForceText = $True }
    @{ Cell = 'E47'; Value = '  -0.34%  '; ForceText = $False }
    @{ Cell = 'D48'; Value = '1.48'; ForceText = $True }
    @{ Cell = 'E48'; Value = '  +1.67%  '; ForceText = $False }
    @{ Cell = 'E49'; Value = '  -0.23%  '; ForceText = $False }
    @{ Cell = 'D50'; Value = '7.55'; ForceText = $True }
    @{ Cell = 'E50'; Value = '  -2.66%  '; ForceText = $False }
    @{ Cell = 'E51'; Value = '  +0.15%  '; ForceText = $False }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
    if ($u.ForceText) {
        $rng.Style = "Normal"
    }
}
